$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 356.45
$ws.Range("I2").Value = 339.2857
$ws.Range("K2").Value = 339.2857
$ws.Range("M2").Value = -226.2857
$ws.Range("H42").Value = 189.73334
$ws.Range("I42").Value = 108.666664
$ws.Range("K42").Value = 325.999992
$ws.Range("M42").Value = -95.99999200000002
$ws.Range("H55").Value = 327.4
$ws.Range("J55").Value = 274.75
$ws.Range("L55").Value = 274.75
$ws.Range("N55").Value = -702.75
$ws.Range("H69").Value = 4702.5
$ws.Range("I69").Value = 2000
$ws.Range("K69").Value = 6000
$ws.Range("M69").Value = -5126
$ws.Range("H70").Value = 6543
$ws.Range("I70").Value = 1702
$ws.Range("K70").Value = 5106
$ws.Range("M70").Value = -4836
$ws.Range("H72").Value = 4702.5
$ws.Range("I72").Value = 2000
$ws.Range("K72").Value = 18000
$ws.Range("M72").Value = -13632
$ws.Range("H73").Value = 6543
$ws.Range("I73").Value = 1702
$ws.Range("K73").Value = 5106
$ws.Range("M73").Value = -4170
$ws.Range("H137").Value = 1834.9
$ws.Range("I137").Value = 1778.5714
$ws.Range("K137").Value = 5335.7142
$ws.Range("M137").Value = -2785.7142
$ws.Range("H138").Value = 2949.52
$ws.Range("I138").Value = 1666.9744
$ws.Range("J138").Value = 3769.5083
$ws.Range("K138").Value = 5000.9232
$ws.Range("L138").Value = 11308.5249
$ws.Range("M138").Value = 139.0767999999998
$ws.Range("N138").Value = -21588.5249
$ws.Range("H141").Value = 1742.9546
$ws.Range("I141").Value = 1519.2222
$ws.Range("K141").Value = 4557.6666
$ws.Range("M141").Value = 622.3334000000004

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 1598.32
$ws.Range("I5").Value = 260.92307
$ws.Range("J5").Value = 3047.1667
$ws.Range("K5").Value = 260.92307
$ws.Range("L5").Value = 3047.1667
$ws.Range("M5").Value = -148.92307
$ws.Range("N5").Value = -3271.1667
$ws.Range("H32").Value = 5261
$ws.Range("I32").Value = 4882.4595
$ws.Range("K32").Value = 4882.4595
$ws.Range("M32").Value = -4595.4595
$ws.Range("H37").Value = 36666.668
$ws.Range("J37").Value = 50000
$ws.Range("L37").Value = 50000
$ws.Range("N37").Value = -50546
$ws.Range("H110").Value = 7609.619
$ws.Range("I110").Value = 9207.643
$ws.Range("K110").Value = 9207.643
$ws.Range("M110").Value = -7162.643
$ws.Range("H122").Value = 3657.7144
$ws.Range("I122").Value = 3421.8
$ws.Range("J122").Value = 4247.5
$ws.Range("K122").Value = 10265.4
$ws.Range("L122").Value = 12742.5
$ws.Range("M122").Value = -7815.400000000001
$ws.Range("N122").Value = -17642.5
$ws.Range("H131").Value = 81710.57000000001
$ws.Range("J131").Value = 81710.57000000001
$ws.Range("L131").Value = 81710.57000000001
$ws.Range("N131").Value = -91790.57000000001
$ws.Range("H132").Value = 2699.5532
$ws.Range("I132").Value = 2396.725
$ws.Range("J132").Value = 4430
$ws.Range("K132").Value = 7190.174999999999
$ws.Range("L132").Value = 13290
$ws.Range("M132").Value = -4660.174999999999
$ws.Range("N132").Value = -18350

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 1598.32
$ws.Range("I4").Value = 260.92307
$ws.Range("J4").Value = 3047.1667
$ws.Range("K4").Value = 260.92307
$ws.Range("L4").Value = 3047.1667
$ws.Range("M4").Value = -145.92307
$ws.Range("N4").Value = -3277.1667
$ws.Range("H20").Value = 3563.68
$ws.Range("I20").Value = 3068.3845
$ws.Range("J20").Value = 4100.25
$ws.Range("K20").Value = 3068.3845
$ws.Range("L20").Value = 4100.25
$ws.Range("M20").Value = -2821.3845
$ws.Range("N20").Value = -4594.25
$ws.Range("H105").Value = 2431.4546
$ws.Range("I105").Value = 1390.0769
$ws.Range("K105").Value = 1390.0769
$ws.Range("M105").Value = 356.9231
$ws.Range("H107").Value = 2952.5454
$ws.Range("I107").Value = 2872.75
$ws.Range("K107").Value = 2872.75
$ws.Range("M107").Value = -952.75

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 248.66667
$ws.Range("I7").Value = 262.18182
$ws.Range("K7").Value = 262.18182
$ws.Range("M7").Value = -149.18182
$ws.Range("H31").Value = 4789.077
$ws.Range("I31").Value = 2009.1875
$ws.Range("K31").Value = 2009.1875
$ws.Range("M31").Value = -1714.1875
$ws.Range("H34").Value = 4789.077
$ws.Range("I34").Value = 2009.1875
$ws.Range("K34").Value = 2009.1875
$ws.Range("M34").Value = -1807.1875
$ws.Range("H134").Value = 12403.333
$ws.Range("I134").Value = 5785.4287
$ws.Range("K134").Value = 17356.2861
$ws.Range("M134").Value = -14821.2861

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 49625000
$ws.Range("I9").Value = 99000000
$ws.Range("J9").Value = 250000
$ws.Range("K9").Value = 297000000
$ws.Range("L9").Value = 750000
$ws.Range("M9").Value = -296999776
$ws.Range("N9").Value = -750448
$ws.Range("H16").Value = 199
$ws.Range("I16").Value = 199
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 597
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -424
$ws.Range("N16").ClearContents()
$ws.Range("H22").Value = 3366.1904
$ws.Range("I22").Value = 448.66666
$ws.Range("J22").Value = 3852.4443
$ws.Range("K22").Value = 1345.99998
$ws.Range("L22").Value = 11557.3329
$ws.Range("M22").Value = -1176.99998
$ws.Range("N22").Value = -11895.3329
$ws.Range("H27").Value = 3366.1904
$ws.Range("I27").Value = 448.66666
$ws.Range("J27").Value = 3852.4443
$ws.Range("K27").Value = 1345.99998
$ws.Range("L27").Value = 11557.3329
$ws.Range("M27").Value = -1243.99998
$ws.Range("N27").Value = -11761.3329
$ws.Range("H34").Value = 3969738.5
$ws.Range("J34").Value = 4631346
$ws.Range("L34").Value = 13894038
$ws.Range("N34").Value = -13894206
$ws.Range("H39").Value = 2070.7144
$ws.Range("J39").Value = 2070.7144
$ws.Range("L39").Value = 6212.1432
$ws.Range("N39").Value = -6800.1432
$ws.Range("H64").Value = 2330
$ws.Range("I64").Value = 1995
$ws.Range("K64").Value = 5985
$ws.Range("M64").Value = -5715
$ws.Range("H67").Value = 2330
$ws.Range("I67").Value = 1995
$ws.Range("K67").Value = 5985
$ws.Range("M67").Value = -5049
$ws.Range("H70").Value = 2633.3333
$ws.Range("I70").Value = 2633.3333
$ws.Range("K70").Value = 7899.999899999999
$ws.Range("M70").Value = -7584.999899999999
$ws.Range("H73").Value = 2633.3333
$ws.Range("I73").Value = 2633.3333
$ws.Range("K73").Value = 7899.999899999999
$ws.Range("M73").Value = -6807.999899999999
$ws.Range("H88").Value = 33334524
$ws.Range("J88").Value = 33334524
$ws.Range("L88").Value = 100003572
$ws.Range("N88").Value = -100004428
$ws.Range("H91").Value = 33334524
$ws.Range("J91").Value = 33334524
$ws.Range("L91").Value = 100003572
$ws.Range("N91").Value = -100006536
$ws.Range("H100").Value = 3962.6667
$ws.Range("J100").Value = 3962.6667
$ws.Range("L100").Value = 11888.0001
$ws.Range("N100").Value = -13510.0001
$ws.Range("H103").Value = 5207.3335
$ws.Range("J103").Value = 7749.75
$ws.Range("L103").Value = 23249.25
$ws.Range("N103").Value = -25007.25
$ws.Range("H109").Value = 498.6
$ws.Range("I109").Value = 498.6
$ws.Range("K109").Value = 1495.8
$ws.Range("M109").Value = -455.8000000000002

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 47620330
$ws.Range("J102").Value = 200002500
$ws.Range("L102").Value = 200002500
$ws.Range("N102").Value = -200005744
$ws.Range("H122").Value = 3120.75
$ws.Range("I122").Value = 900
$ws.Range("J122").Value = 3861
$ws.Range("K122").Value = 2700
$ws.Range("L122").Value = 11583
$ws.Range("M122").Value = -250
$ws.Range("N122").Value = -16483
$ws.Range("H126").Value = 24210.732
$ws.Range("I126").Value = 47848.145
$ws.Range("J126").Value = 3528
$ws.Range("K126").Value = 143544.435
$ws.Range("L126").Value = 10584
$ws.Range("M126").Value = -141074.435
$ws.Range("N126").Value = -15524
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4930.5
$ws.Range("I40").Value = 3733.4119
$ws.Range("K40").Value = 3733.4119
$ws.Range("M40").Value = -3597.4119
$ws.Range("H61").Value = 4480.4707
$ws.Range("I61").Value = 4416.8125
$ws.Range("K61").Value = 4416.8125
$ws.Range("M61").Value = -4214.8125
$ws.Range("H93").Value = 2215
$ws.Range("I93").Value = 2215
$ws.Range("K93").Value = 2215
$ws.Range("M93").Value = -967
$ws.Range("H113").Value = 4480.4707
$ws.Range("I113").Value = 4416.8125
$ws.Range("K113").Value = 4416.8125
$ws.Range("M113").Value = -2246.8125

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 476.44446
$ws.Range("I107").Value = 476.44446
$ws.Range("K107").Value = 1429.33338
$ws.Range("M107").Value = 490.66662

